$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("G2").Value = "2016-08-25 15:10:27"
$zhcn.Range("H2").Value = "2016-08-25 15:10:17"
$zhcn.Range("K2").Value = "2016-08-25 15:10:47"
$dede.Range("H2").Value = "2016-08-25 15:10:27"
$dede.Range("K2").Value = "2016-08-25 15:10:54"
